$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - Cantar (scale)
$ws.Range("A20").Value = "Cantar "
$ws.Range("A20").WrapText = $true
$ws.Range("B20").Value = "https://www.galaxus.ch/im/productimages/1/5/9/1/5/3/9/8/8/7/2/6/3/7/8/1/9/0/6/bf70da80-90ba-4ff6-a7e0-7da33a3eb186_cropped.jpg?impolicy=product&resizeWidth=720"
$ws.Range("D20").Value = "80 - 120 CHF"

# Row 21 - Papuci 40 (slippers)
$ws.Range("A21").Value = "Papuci 40"
$ws.Range("A21").WrapText = $true
$ws.Range("B21").Value = "https://img01.ztat.net/article/spp-media-p1/fa20de609217474a85fe92b8da997df0/8a35b97c476f4d8aa1b9765d103a51b3.jpg?imwidth=1800&filter=packshot"
$ws.Range("C21").Value = "https://fr.zalando.ch/anna-field-chaussons-grey-an611d024-c11.html?_rfl=de"
$ws.Range("D21").Value = "34 CHF"

# Row 22 - Lingura (spoon)
$ws.Range("A22").Value = "Lingura"
$ws.Range("A22").WrapText = $true
$ws.Range("B22").Value = "https://m.media-amazon.com/images/I/71bxju0FpgL._AC_SL1200_.jpg"
$ws.Range("C22").Value = "https://www.amazon.de/dp/B01EZAC8P0?ref_=cm_sw_r_apin_ct_12HQWTDZT7YJGK9MFAPB_1&language=en_US&th=1"
$ws.Range("D22").Value = "20 EUR"

# Row 23 - Set linguri (spoon set)
$ws.Range("A23").Value = "Set linguri"
$ws.Range("A23").WrapText = $true
$ws.Range("B23").Value = "https://m.media-amazon.com/images/I/61g0EHYpeWL._AC_SL1500_.jpg"
$ws.Range("C23").Value = "https://www.amazon.de/dp/B00VRK60ZQ?ref_=cm_sw_r_apin_ct_12HQWTDZT7YJGK9MFAPB&language=en_US&th=1"
$ws.Range("D23").Value = "20 EUR"

# Row 24 - Pensule (brushes)
$ws.Range("A24").Value = "Pensule"
$ws.Range("A24").WrapText = $true
$ws.Range("B24").Value = "https://m.media-amazon.com/images/I/615gJtEdGFL._AC_SL1500_.jpg"
$ws.Range("C24").Value = "https://www.amazon.de/dp/B09C29PHSW?ref=cm_sw_r_cso_wa_apin_dp_4MN968FMNN3XZE83RGKE&ref_=cm_sw_r_cso_wa_apin_dp_4MN968FMNN3XZE83RGKE&social_share=cm_sw_r_cso_wa_apin_dp_4MN968FMNN3XZE83RGKE&starsLeft=1&skipTwisterOG=1"
$ws.Range("D24").Value = "6 EUR"

# Row 25 - Fata perna (pillowcase)
$ws.Range("A25").Value = "Fata perna"
$ws.Range("A25").WrapText = $true
$ws.Range("B25").Value = "https://m.media-amazon.com/images/I/61dBnRJ8qPL._AC_SL1500_.jpg"
$ws.Range("C25").Value = "https://www.amazon.de/dp/B08B4BM4DW?ref=cm_sw_r_cso_wa_apin_dp_Q9MP39GHH5HJTFKPPC1V&ref_=cm_sw_r_cso_wa_apin_dp_Q9MP39GHH5HJTFKPPC1V&social_share=cm_sw_r_cso_wa_apin_dp_Q9MP39GHH5HJTFKPPC1V&starsLeft=1&skipTwisterOG=1&th=1"
$ws.Range("D25").Value = "29 EUR"

# Update selection to match the target state (B22 selected)
$ws.Range("B22").Select() | Out-Null
